$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Programs" SQL query in B2 to add the CASE expression for the
# "Website" column (program_link / program_acronym logic).
$newProgramsQuery = @"
SELECT DISTINCT 
    prg.program_name AS "Program",
  CASE
    WHEN prg.program_link IS NOT NULL THEN prg.program_acronym
        ELSE prg.program_link
    END  AS "Website",
    prg.focus_area AS "Focus Area",
    prg.cancer_type AS "Cancer Type",
 CASE 
        WHEN prg.data_link IS NOT NULL THEN prg.website       
        ELSE prg.data_link
    END AS "Data Location Details"
FROM 
    df_program prg
WHERE 
     prg.cancer_type LIKE '%Ovarian Cancer%'
ORDER BY 
    lower(prg.program_name) ASC
LIMIT 100;
"@

$ws.Range("B2").Value = $newProgramsQuery

# Move the current selection from C2 to C3 and drop the scrolled
# topLeftCell state (scroll back to show row 1).
$ws.Range("C3").Select() | Out-Null
